$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "ROW50-FE-LIFTER": append row 56 (new reading, same shape as row 55)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
$r1 = 56
$ws1.Cells.Item($r1, 1).NumberFormat = $ws1.Cells.Item($r1 - 1, 1).NumberFormat
$ws1.Cells.Item($r1, 1).Value = 45752.20712444445
$ws1.Cells.Item($r1, 2).Value = "0x01,0x90"
$ws1.Cells.Item($r1, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws1.Cells.Item($r1, 4).Value = "0x01,0x5e"
$ws1.Cells.Item($r1, 5).Value = "0xe"
$ws1.Cells.Item($r1, 6).Value = 400
$ws1.Cells.Item($r1, 7).Value = [double]"5.68631262647114e+23"
$ws1.Cells.Item($r1, 8).Value = 350
$ws1.Cells.Item($r1, 9).Value = 14

# ---------------------------------------------------------------------------
# Sheet "ROW50-MID-LIFTER": append row 58
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
$r2 = 58
$ws2.Cells.Item($r2, 1).NumberFormat = $ws2.Cells.Item($r2 - 1, 1).NumberFormat
$ws2.Cells.Item($r2, 1).Value = 45752.17625
$ws2.Cells.Item($r2, 2).Value = "0x01,0x90 "
$ws2.Cells.Item($r2, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Cells.Item($r2, 4).Value = "0x01,0x62"
$ws2.Cells.Item($r2, 5).Value = "0x19"
$ws2.Cells.Item($r2, 6).Value = 400
# This particular value is too large to round-trip as a clean double, so the
# source workbook keeps it as literal text instead of a number.
$ws2.Cells.Item($r2, 7).NumberFormat = "@"
$ws2.Cells.Item($r2, 7).Value = "568631262647113771663628"
$ws2.Cells.Item($r2, 7).Style = "Normal"
$ws2.Cells.Item($r2, 8).Value = 354
$ws2.Cells.Item($r2, 9).Value = 25

# ---------------------------------------------------------------------------
# Sheet "ROW11-FE-LIFTER": append row 56
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
$r3 = 56
$ws3.Cells.Item($r3, 1).NumberFormat = $ws3.Cells.Item($r3 - 1, 1).NumberFormat
$ws3.Cells.Item($r3, 1).Value = 45752.23778543981
$ws3.Cells.Item($r3, 2).Value = "0x01,0x90"
$ws3.Cells.Item($r3, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws3.Cells.Item($r3, 4).Value = "0x01,0x5e"
$ws3.Cells.Item($r3, 5).Value = "0x14"
$ws3.Cells.Item($r3, 6).Value = 400
$ws3.Cells.Item($r3, 7).Value = [double]"5.68631262647114e+23"
$ws3.Cells.Item($r3, 8).Value = 350
$ws3.Cells.Item($r3, 9).Value = 20

# ---------------------------------------------------------------------------
# Sheet "ROW11-MID-LIFTER": append row 56
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
$r4 = 56
$ws4.Cells.Item($r4, 1).NumberFormat = $ws4.Cells.Item($r4 - 1, 1).NumberFormat
$ws4.Cells.Item($r4, 1).Value = 45752.3739091088
$ws4.Cells.Item($r4, 2).Value = "0x01,0x90"
$ws4.Cells.Item($r4, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws4.Cells.Item($r4, 4).Value = "0x01,0x66"
$ws4.Cells.Item($r4, 5).Value = "0x19"
$ws4.Cells.Item($r4, 6).Value = 400
$ws4.Cells.Item($r4, 7).Value = [double]"5.68631262647114e+23"
$ws4.Cells.Item($r4, 8).Value = 358
$ws4.Cells.Item($r4, 9).Value = 25
